$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all modified cells so numeric-looking strings
# (e.g. "0.5250", "220.01") are preserved verbatim as text, matching
# the source data which stores every Coin/Link/Price/Volume value as a string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.388.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5250"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2670"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06357"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07763"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.669.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.463"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5516"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8264"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.400.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.24"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.26"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.257"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1263"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.96"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.374"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.18"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.417"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06103"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.291"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.581"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.393"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.675"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9983"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.423"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6029"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.781"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01608"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.993"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.083.39"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8556"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.811.10"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.78"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.115"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05204"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.475"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4233"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.59%  "
